$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 0.01118366666666667
$ws.Range("N2").Value = 0.033551
$ws.Range("O2").Value = 0.1526349789820392
$ws.Range("P2").Value = 0.1526349789820392
$ws.Range("Q2").Value = 0.001610313796
$ws.Range("R2").Value = 0.014492824164
$ws.Range("S2").Value = 0.1526349789820392
$ws.Range("T2").Value = 0.1526349789820392

# Row 3
$ws.Range("O3").Value = 0.3287218168252871
$ws.Range("P3").Value = 0.3287218168252871
$ws.Range("S3").Value = 0.3287218168252871
$ws.Range("T3").Value = 0.3287218168252871

# Row 4
$ws.Range("M4").Value = 0.03800133333333333
$ws.Range("N4").Value = 0.114004
$ws.Range("O4").Value = 0.5186432041926737
$ws.Range("P4").Value = 0.5186432041926737
$ws.Range("Q4").Value = 0.005471735984
$ws.Range("R4").Value = 0.049245623856
$ws.Range("S4").Value = 0.5186432041926737
$ws.Range("T4").Value = 0.5186432041926737
